# export pdf control version
# Applies three localized edits to the document:
#  1. First (empty) paragraph mark: add w:hint="cs" + <w:cs/> to its rPr.
#  2. The "table caption" paragraph: "ตาราง" + " " + "…" + " " + " " + "State Machine Diagram"
#     becomes "ตาราง" + "ที่ " + "1 " + "State Machine Diagram".
#  3. The "กิตติพศ " run in the version-control table is split into
#     "กิตติ" + <proofErr spellStart> + "พศ" + <proofErr spellEnd> + " ".
#
# Each edit is performed by rebuilding the *entire* containing paragraph's
# OOXML (preserving every original attribute untouched by the diff, e.g.
# w14:paraId / w:rsidR / w:rsidRPr on the paragraph and on every run that
# is not itself being changed) and pushing it back in with Range.InsertXML,
# which is the only reliable way to get <w:cs/> / <w:proofErr/> markup
# exactly right (these aren't modeled as ordinary Word object-model
# properties).

$d = $word.ActiveDocument

function New-OoxmlPackage([string]$bodyXml) {
  return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- Edit 1: paragraph-mark formatting on the first (empty) paragraph ----
$para1 = $d.Paragraphs(1)
$p1Xml = '<w:p w14:paraId="02A8AA18" w14:textId="77777777" w:rsidR="00667B81" w:rsidRDefault="00667B81" w:rsidP="00667B81">' +
  '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:cs/></w:rPr></w:pPr>' +
  '</w:p>'
$para1.Range.InsertXML((New-OoxmlPackage $p1Xml))

# --- Edit 2: "ตาราง ... State Machine Diagram" caption paragraph ---------
$rng2 = $d.Content.Duplicate
[void]$rng2.Find.Execute("ตาราง", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1)
$p2Xml = '<w:p w14:paraId="7B875441" w14:textId="4E0CD68B" w:rsidR="00667B81" w:rsidRDefault="00667B81" w:rsidP="00667B81">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00C60201"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t>ตาราง</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve">ที่ </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">1 </w:t></w:r>' +
  '<w:r w:rsidR="00EB3896" w:rsidRPr="00EB3896"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">State </w:t></w:r>' +
  '<w:r w:rsidR="00EB3896"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>M</w:t></w:r>' +
  '<w:r w:rsidR="00EB3896" w:rsidRPr="00EB3896"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">achine </w:t></w:r>' +
  '<w:r w:rsidR="00EB3896"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>D</w:t></w:r>' +
  '<w:r w:rsidR="00EB3896" w:rsidRPr="00EB3896"><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>iagram</w:t></w:r>' +
  '</w:p>'
$para2.Range.InsertXML((New-OoxmlPackage $p2Xml))

# --- Edit 3: split the "กิตติพศ " run and wrap "พศ" with proofErr tags ---
$rng3 = $d.Content.Duplicate
[void]$rng3.Find.Execute("กิตติพศ (SP)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para3 = $rng3.Paragraphs(1)
$p3Xml = '<w:p w14:paraId="63296871" w14:textId="77777777" w:rsidR="00667B81" w:rsidRDefault="00667B81" w:rsidP="00530DA5">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t>กิตติ</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t>พศ</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New" w:hint="cs"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>(SP)</w:t></w:r>' +
  '</w:p>'
$para3.Range.InsertXML((New-OoxmlPackage $p3Xml))

Write-Output "edits applied"
